$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.049.42'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.09%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.457.39'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.44%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.54%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '179.21'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.60%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.609'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.53%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.457.22'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.39%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.137'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.84%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.94'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.03%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.431'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.29%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.060.04'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.12%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '31.65'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +6.35%  '

# Row 15
$ws.Range('E15').Value = '  -0.29%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.023.20'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.09%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000176'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.60%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.464.54'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.06%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.24'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.48%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.12'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.07%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '388.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.85%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.90'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.49%  '

# Row 23
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.12%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.38'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.43%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.74'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.99%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.535'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.13%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000122'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.16%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.33'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.59%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.174'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.23%  '

# Row 30
$ws.Range('E30').Value = '  +0.59%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.15'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.29%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.40'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.47%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.05'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.29%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.43'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.84%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.34'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.31%  '

# Row 36
$ws.Range('E36').Value = '  -0.04%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.59'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.13%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '161.92'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.01%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.877'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.02%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.80'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +8.90%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.86'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.66%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.80'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.68%  '

# Row 43
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.40%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '26.05'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.29%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.791.48'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.87%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0719'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.64%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.29'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.11%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '41.05'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.33%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0297'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.78%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '328.16'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.75%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.04'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.39%  '
